$wb = $excel.ActiveWorkbook
$ws6 = $wb.Worksheets.Item("list 6")

# --- "list 6" (dog/cat breed list): rename the Yorkshire Terrier entry to the
# new cross-breed name, and add a brand-new "Mixed-breed dog" / "Croisé" row
# right under it (shifting the old "Rubber Duck" row down by one). ---

# Row 11 used to read "Yorkshire Terrier" / "Terrier du Yorkshire"; it now
# reads "Yorkshire Terrier-Maltese Mix" / "Terrier du Yorkshire-Maltese Mix".
# The English-column cell also picks up an explicit black font color.
$ws6.Cells.Item(11, 6).Value = "Yorkshire Terrier-Maltese Mix"
$ws6.Cells.Item(11, 6).Font.Color = 0
$ws6.Cells.Item(11, 7).Value = "Terrier du Yorkshire-Maltese Mix"

# Insert a new row 12 (inherits formatting from row 11) for the new breed
# entry, pushing the former row 12 ("Rubber Duck") down to row 13.
$ws6.Rows.Item(12).Insert()
$ws6.Cells.Item(12, 2).Value = "Specie"
$ws6.Cells.Item(12, 3).Value = "Espèce"
$ws6.Cells.Item(12, 4).Value = "Dog"
$ws6.Cells.Item(12, 5).Value = "Chien"
$ws6.Cells.Item(12, 7).Value = "Croisé"
$ws6.Cells.Item(12, 6).Value = "Mixed-breed dog"
$ws6.Cells.Item(12, 6).Font.Color = 0

# New row 12's license cell (column D) gets the same CC BY-SA hyperlink the
# other breed rows already carry.
[void]$ws6.Hyperlinks.Add($ws6.Cells.Item(12, 4), "https://creativecommons.org/licenses/by-sa/4.0/", "", "", "https://creativecommons.org/licenses/by-sa/4.0/")

# Widen columns F and G so the longer new breed names fit.
$ws6.Columns.Item(6).EntireColumn.AutoFit()
$ws6.Columns.Item(7).EntireColumn.AutoFit()

# The former blank row 13 (now row 14) loses its stray D-column cell.
$ws6.Cells.Item(14, 4).ClearContents()

# "list 6" becomes the active/visible sheet (it was "list 5" before).
[void]$ws6.Activate()
[void]$ws6.Select()
[void]$ws6.Range("E9").Select()
